$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.113.90"
$ws.Range("E2").Value = "  -3.94%  "

$ws.Range("D3").Value = "3.523.38"
$ws.Range("E3").Value = "  -3.19%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'578.21"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").Value = "'172.18"
$ws.Range("E6").Value = "  -5.26%  "

$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  -0.52%  "

$ws.Range("D8").Value = "3.515.57"
$ws.Range("E8").Value = "  -3.14%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'0.189"
$ws.Range("E10").Value = "  -7.04%  "

$ws.Range("D11").Value = "'6.71"
$ws.Range("E11").Value = "  +11.28%  "

$ws.Range("D12").Value = "'0.605"
$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").Value = "'47.29"
$ws.Range("E13").Value = "  -5.14%  "

$ws.Range("E14").Value = "  -4.11%  "

$ws.Range("D15").Value = "'692.54"
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("D16").Value = "4.091.83"
$ws.Range("E16").Value = "  -3.10%  "

$ws.Range("D17").Value = "'8.85"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "69.178.15"
$ws.Range("E18").Value = "  -3.84%  "

$ws.Range("D19").Value = "3.526.27"
$ws.Range("E19").Value = "  -2.92%  "

$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("D21").Value = "'17.48"
$ws.Range("E21").Value = "  -4.59%  "

$ws.Range("D22").Value = "'11.18"
$ws.Range("E22").Value = "  -4.05%  "

$ws.Range("D23").Value = "'0.906"
$ws.Range("E23").Value = "  -3.91%  "

$ws.Range("D24").Value = "'16.60"
$ws.Range("E24").Value = "  -6.75%  "

$ws.Range("D25").Value = "'97.73"
$ws.Range("E25").Value = "  -5.37%  "

$ws.Range("E26").Value = "  -4.42%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -6.51%  "

$ws.Range("E29").Value = "  -6.03%  "

$ws.Range("D30").Value = "'33.24"
$ws.Range("E30").Value = "  -5.95%  "

$ws.Range("D31").Value = "'8.87"
$ws.Range("E31").Value = "  -4.16%  "

$ws.Range("E32").Value = "  -7.05%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.29"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "'1.35"
$ws.Range("E34").Value = "  -5.96%  "

$ws.Range("D35").Value = "'564.54"
$ws.Range("E35").Value = "  -3.62%  "

$ws.Range("D36").Value = "'3.65"
$ws.Range("E36").Value = "  -12.79%  "

$ws.Range("D37").Value = "'10.86"
$ws.Range("E37").Value = "  -4.17%  "

$ws.Range("E38").Value = "  -3.38%  "

$ws.Range("D39").Value = "'57.27"
$ws.Range("E39").Value = "  -3.92%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0441"
$ws.Range("E41").Value = "  -6.41%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.339"
$ws.Range("E42").Value = "  -2.45%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.138"
$ws.Range("E43").Value = "  -3.92%  "

$ws.Range("D44").Value = "3.447.38"
$ws.Range("E44").Value = "  -6.97%  "

$ws.Range("D45").Value = "'33.26"
$ws.Range("E45").Value = "  -7.07%  "

$ws.Range("E46").Value = "  -8.24%  "

$ws.Range("E47").Value = "  +2.11%  "

$ws.Range("E48").Value = "  -7.63%  "

$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("D50").Value = "'134.35"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("E51").Value = "  -0.37%  "
